# Generate Report for Handoff
# Update the localization-status report: the zh-cn / de-de languages have
# moved from "In Translation" to "Ready for handoff", and the corresponding
# timestamps have been refreshed.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-08-20 18:49:56"

# Columns E and F grew wider to accommodate the new, longer status text.
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-08-20 18:49:52"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333332

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-08-20 18:49:56"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333332
